$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 29.705
$ws.Cells.Item(2, 5).Value = 5.624
$ws.Cells.Item(2, 6).Value = 5.051
$ws.Cells.Item(2, 7).Value = 1.019
$ws.Cells.Item(2, 10).Value = 29.705
$ws.Cells.Item(2, 11).Value = 5.624
$ws.Cells.Item(2, 13).Value = 0.628
$ws.Cells.Item(2, 16).Value = 118.423
$ws.Cells.Item(2, 17).Value = 5.624
$ws.Cells.Item(2, 18).Value = 5.051
$ws.Cells.Item(2, 19).Value = 0.504
$ws.Cells.Item(3, 4).Value = 29.705
$ws.Cells.Item(3, 5).Value = 5.624
$ws.Cells.Item(3, 6).Value = 5.051
$ws.Cells.Item(3, 7).Value = 1.019
$ws.Cells.Item(3, 10).Value = 29.705
$ws.Cells.Item(3, 11).Value = 5.624
$ws.Cells.Item(3, 13).Value = 0.628
$ws.Cells.Item(3, 16).Value = 118.423
$ws.Cells.Item(3, 17).Value = 5.624
$ws.Cells.Item(3, 18).Value = 5.051
$ws.Cells.Item(3, 19).Value = 0.504
$ws.Cells.Item(4, 2).Value = -93.235
$ws.Cells.Item(4, 5).Value = 5.624
$ws.Cells.Item(4, 6).Value = 4.831
$ws.Cells.Item(4, 11).Value = 5.624
$ws.Cells.Item(4, 12).Value = 3.228
$ws.Cells.Item(4, 13).Value = 0.533
$ws.Cells.Item(4, 14).Value = -30.974
$ws.Cells.Item(4, 16).Value = 118.044
$ws.Cells.Item(4, 17).Value = 5.624
$ws.Cells.Item(4, 18).Value = 4.831
$ws.Cells.Item(4, 19).Value = 0.481
$ws.Cells.Item(5, 2).Value = -93.371
$ws.Cells.Item(5, 4).Value = 29.612
$ws.Cells.Item(5, 5).Value = 5.624
$ws.Cells.Item(5, 7).Value = 0.769
$ws.Cells.Item(5, 8).Value = -93.89700000000001
$ws.Cells.Item(5, 10).Value = 29.612
$ws.Cells.Item(5, 11).Value = 5.624
$ws.Cells.Item(5, 12).Value = 3.239
$ws.Cells.Item(5, 13).Value = 0.527
$ws.Cells.Item(5, 14).Value = -31.019
$ws.Cells.Item(5, 16).Value = 118.049
$ws.Cells.Item(5, 17).Value = 5.624
$ws.Cells.Item(5, 19).Value = 0.481
$ws.Cells.Item(6, 2).Value = -92.718
$ws.Cells.Item(6, 4).Value = 30.064
$ws.Cells.Item(6, 5).Value = 5.624
$ws.Cells.Item(6, 6).Value = 4.836
$ws.Cells.Item(6, 7).Value = 0.773
$ws.Cells.Item(6, 8).Value = -93.708
$ws.Cells.Item(6, 10).Value = 30.064
$ws.Cells.Item(6, 11).Value = 5.624
$ws.Cells.Item(6, 12).Value = 3.237
$ws.Cells.Item(6, 13).Value = 0.527
$ws.Cells.Item(6, 14).Value = -30.802
$ws.Cells.Item(6, 16).Value = 119.854
$ws.Cells.Item(6, 17).Value = 5.624
$ws.Cells.Item(6, 18).Value = 4.836
$ws.Cells.Item(7, 2).Value = -92.224
$ws.Cells.Item(7, 4).Value = 52.547
$ws.Cells.Item(7, 5).Value = 5.624
$ws.Cells.Item(7, 6).Value = 4.842
$ws.Cells.Item(7, 7).Value = 0.662
$ws.Cells.Item(7, 8).Value = -93.258
$ws.Cells.Item(7, 10).Value = 52.547
$ws.Cells.Item(7, 11).Value = 5.624
$ws.Cells.Item(7, 12).Value = 3.234
$ws.Cells.Item(7, 13).Value = 0.473
$ws.Cells.Item(7, 14).Value = -30.638
$ws.Cells.Item(7, 16).Value = 209.481
$ws.Cells.Item(7, 17).Value = 5.624
$ws.Cells.Item(7, 18).Value = 4.842
$ws.Cells.Item(7, 19).Value = 0.176
$ws.Cells.Item(8, 2).Value = -92.224
$ws.Cells.Item(8, 4).Value = 52.547
$ws.Cells.Item(8, 5).Value = 5.624
$ws.Cells.Item(8, 6).Value = 4.842
$ws.Cells.Item(8, 7).Value = 0.662
$ws.Cells.Item(8, 8).Value = -93.258
$ws.Cells.Item(8, 10).Value = 52.547
$ws.Cells.Item(8, 11).Value = 5.624
$ws.Cells.Item(8, 12).Value = 3.234
$ws.Cells.Item(8, 13).Value = 0.473
$ws.Cells.Item(8, 14).Value = -30.638
$ws.Cells.Item(8, 16).Value = 209.481
$ws.Cells.Item(8, 17).Value = 5.624
$ws.Cells.Item(8, 18).Value = 4.842
$ws.Cells.Item(8, 19).Value = 0.176
$ws.Cells.Item(9, 2).Value = -92.718
$ws.Cells.Item(9, 4).Value = 30.064
$ws.Cells.Item(9, 5).Value = 5.624
$ws.Cells.Item(9, 6).Value = 4.836
$ws.Cells.Item(9, 7).Value = 0.773
$ws.Cells.Item(9, 8).Value = -93.708
$ws.Cells.Item(9, 10).Value = 30.064
$ws.Cells.Item(9, 11).Value = 5.624
$ws.Cells.Item(9, 12).Value = 3.237
$ws.Cells.Item(9, 13).Value = 0.527
$ws.Cells.Item(9, 14).Value = -30.802
$ws.Cells.Item(9, 16).Value = 119.854
$ws.Cells.Item(9, 17).Value = 5.624
$ws.Cells.Item(9, 18).Value = 4.836
$ws.Cells.Item(10, 2).Value = -93.371
$ws.Cells.Item(10, 4).Value = 29.612
$ws.Cells.Item(10, 5).Value = 5.624
$ws.Cells.Item(10, 7).Value = 0.769
$ws.Cells.Item(10, 8).Value = -93.89700000000001
$ws.Cells.Item(10, 10).Value = 29.612
$ws.Cells.Item(10, 11).Value = 5.624
$ws.Cells.Item(10, 12).Value = 3.239
$ws.Cells.Item(10, 13).Value = 0.527
$ws.Cells.Item(10, 14).Value = -31.019
$ws.Cells.Item(10, 16).Value = 118.049
$ws.Cells.Item(10, 17).Value = 5.624
$ws.Cells.Item(10, 19).Value = 0.481
$ws.Cells.Item(11, 2).Value = -93.235
$ws.Cells.Item(11, 5).Value = 5.624
$ws.Cells.Item(11, 6).Value = 4.831
$ws.Cells.Item(11, 11).Value = 5.624
$ws.Cells.Item(11, 12).Value = 3.228
$ws.Cells.Item(11, 13).Value = 0.533
$ws.Cells.Item(11, 14).Value = -30.974
$ws.Cells.Item(11, 16).Value = 118.044
$ws.Cells.Item(11, 17).Value = 5.624
$ws.Cells.Item(11, 18).Value = 4.831
$ws.Cells.Item(11, 19).Value = 0.481
